# Edit: insert 8 new Durazno price rows (Carson, Loadel, Toscana) at row 712
# on the "Femacal de La Calera" sheet, shifting existing rows 712:755 down to 720:763.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 blank rows at 712, pushing existing data down (matches dimension A1:T763 afterwards)
$ws.Rows("712:719").Insert()

$newRows = @(
    @{ D = 44578; K = "Carson"; L = "Especial"; M = 80; N = 16000; O = 16000; P = 16000; Q = "`$/caja 15 kilos empedrada"; R = "Provincia de San Felipe de Aconcagua"; S = 1067; T = 15 },
    @{ D = 44578; K = "Carson"; L = "Primera"; M = 90; N = 14000; O = 14000; P = 14000; Q = "`$/caja 15 kilos empedrada"; R = "Provincia de San Felipe de Aconcagua"; S = 933; T = 15 },
    @{ D = 44578; K = "Carson"; L = "Segunda"; M = 80; N = 12000; O = 12000; P = 12000; Q = "`$/caja 15 kilos empedrada"; R = "Provincia de San Felipe de Aconcagua"; S = 800; T = 15 },
    @{ D = 44578; K = "Loadel"; L = "Especial"; M = 80; N = 15000; O = 15000; P = 15000; Q = "`$/caja 15 kilos empedrada"; R = "Provincia de San Felipe de Aconcagua"; S = 1000; T = 15 },
    @{ D = 44578; K = "Loadel"; L = "Primera"; M = 75; N = 12000; O = 12000; P = 12000; Q = "`$/caja 15 kilos empedrada"; R = "Provincia de San Felipe de Aconcagua"; S = 800; T = 15 },
    @{ D = 44578; K = "Loadel"; L = "Segunda"; M = 70; N = 10000; O = 10000; P = 10000; Q = "`$/caja 15 kilos empedrada"; R = "Provincia de San Felipe de Aconcagua"; S = 667; T = 15 },
    @{ D = 44578; K = "Toscana"; L = "Primera"; M = 80; N = 14000; O = 14000; P = 14000; Q = "`$/caja 15 kilos empedrada"; R = "Región de O'Higgins"; S = 933; T = 15 },
    @{ D = 44578; K = "Toscana"; L = "Segunda"; M = 70; N = 12000; O = 12000; P = 12000; Q = "`$/caja 15 kilos empedrada"; R = "Región de O'Higgins"; S = 800; T = 15 }
)

$r = 712
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = 5
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100103
    $ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value = 100103004
    $ws.Cells.Item($r, 10).Value = "Durazno"
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r = $r + 1
}

Write-Host "Done. Dimension should now be A1:T763."
